# Updates cryptos list values per the Jan 26 2024 GitHub Actions data refresh.
# Column D ("Price") cells are forced to Text format before assignment so that
# numeric-looking price strings (e.g. "87.90", "0.470") keep their exact text
# representation instead of being auto-converted to numbers by Excel.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "40.192.00"
$ws.Range("E2").Value = "  +0.50%  "
# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.226.68"
$ws.Range("E3").Value = "  +0.64%  "
# Row 4
$ws.Range("E4").Value = "  +0.01%  "
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "294.04"
$ws.Range("E5").Value = "  +1.80%  "
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "87.90"
$ws.Range("E6").Value = "  +0.00%  "
# Row 7
$ws.Range("E7").Value = "  -0.28%  "
# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.470"
$ws.Range("E9").Value = "  -0.06%  "
# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "30.84"
$ws.Range("E10").Value = "  +0.55%  "
# Row 11
$ws.Range("E11").Value = "  +6.58%  "
# Row 12
$ws.Range("E12").Value = "  +0.09%  "
# Row 13
$ws.Range("E13").Value = "  +3.18%  "
# Row 14
$ws.Range("E14").Value = "  +0.12%  "
# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.585.87"
$ws.Range("E15").Value = "  +1.21%  "
# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "13.89"
$ws.Range("E16").Value = "  -0.92%  "
# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.233.52"
$ws.Range("E17").Value = "  +2.07%  "
# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.738"
# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "40.118.74"
# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0890"
$ws.Range("E20").Value = "  +0.58%  "
# Row 21
$ws.Range("E21").Value = "  -4.78%  "
# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.79"
$ws.Range("E22").Value = "  -0.04%  "
# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "65.75"
$ws.Range("E23").Value = "  +0.14%  "
# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "236.88"
$ws.Range("E24").Value = "  +0.77%  "
# Row 25
$ws.Range("E25").Value = "  +0.03%  "
# Row 26
$ws.Range("E26").Value = "  +0.69%  "
# Row 27
$ws.Range("E27").Value = "  -0.52%  "
# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "23.32"
$ws.Range("E28").Value = "  +3.27%  "
# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.32"
# Row 30
$ws.Range("E30").Value = "  -6.69%  "
# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "158.10"
$ws.Range("E31").Value = "  +2.79%  "
# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "31.92"
$ws.Range("E32").Value = "  -0.69%  "
# Row 33
$ws.Range("E33").Value = "  -0.07%  "
# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.98"
$ws.Range("E34").Value = "  +0.74%  "
# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.05"
$ws.Range("E35").Value = "  +7.26%  "
# Row 36
$ws.Range("E36").Value = "  -0.20%  "
# Row 37
$ws.Range("E37").Value = "  -2.66%  "
# Row 38
$ws.Range("E38").Value = "  +1.40%  "
# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.76"
$ws.Range("E39").Value = "  +3.05%  "
# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0995"
$ws.Range("E40").Value = "  -0.52%  "
# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "15.51"
$ws.Range("E41").Value = "  -3.09%  "
# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.083.99"
$ws.Range("E42").Value = "  -0.15%  "
# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.74"
$ws.Range("E43").Value = "  -2.11%  "
# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "18.90"
$ws.Range("E44").Value = "  +6.80%  "
# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.14"
$ws.Range("E45").Value = "  +2.79%  "
# Row 46
$ws.Range("E46").Value = "  +0.91%  "
# Row 47
$ws.Range("B47").Value = "NEARProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.75"
$ws.Range("E47").Value = "  +3.10%  "
# Row 48
$ws.Range("B48").Value = "ApeXProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.95"
$ws.Range("E48").Value = "  -11.43%  "
# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.452.53"
$ws.Range("E49").Value = "  +1.09%  "
# Row 50
$ws.Range("E50").Value = "  +3.27%  "
# Row 51
$ws.Range("E51").Value = "  +3.90%  "

Write-Output "Applied 81 cell updates to Sheet1"
